$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray C2 value entirely (cell becomes empty)
$ws.Range("C2").ClearContents()

# Correct floating point precision drifts from the naive forecaster bug fix
$ws.Range("C3").Value = -8.992252553594259
$ws.Range("E3").Value = -19.76480035196673

$ws.Range("C4").Value = 7.007132997505217

$ws.Range("C5").Value = 8.866443976147087

$ws.Range("C6").Value = 3.0013062146236

$ws.Range("C7").Value = -2.90476933598719

$ws.Range("E8").Value = 9.131012060398703

$ws.Range("C11").Value = 4.073887526082043

$ws.Range("E12").Value = -4.308894244053663

$ws.Range("C14").Value = -2.305533699949858

$ws.Range("E15").Value = 16.14645080511215

$ws.Range("E17").Value = -3.570724939213787
